$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1829.7693
$ws.Range("I86").Value = 1818.4546
$ws.Range("J86").Value = 1892
$ws.Range("K86").Value = 1818.4546
$ws.Range("L86").Value = 1892
$ws.Range("M86").Value = -695.4546
$ws.Range("N86").Value = -4138

$ws.Range("H89").Value = 1829.7693
$ws.Range("I89").Value = 1818.4546
$ws.Range("J89").Value = 1892
$ws.Range("K89").Value = 9092.273000000001
$ws.Range("L89").Value = 9460
$ws.Range("M89").Value = -3476.273000000001
$ws.Range("N89").Value = -20692

$ws.Range("H137").Value = 8930510
$ws.Range("J137").Value = 2399.7856
$ws.Range("L137").Value = 7199.3568
$ws.Range("N137").Value = -12299.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25646788
$ws.Range("I32").Value = 32261132
$ws.Range("J32").Value = 16201.75
$ws.Range("K32").Value = 32261132
$ws.Range("L32").Value = 16201.75
$ws.Range("M32").Value = -32260845
$ws.Range("N32").Value = -16775.75

$ws.Range("H88").Value = 1950.8334
$ws.Range("I88").Value = 2015.5
$ws.Range("J88").Value = 1627.5
$ws.Range("K88").Value = 2015.5
$ws.Range("L88").Value = 1627.5
$ws.Range("M88").Value = -1609.5
$ws.Range("N88").Value = -2439.5

$ws.Range("H91").Value = 1950.8334
$ws.Range("I91").Value = 2015.5
$ws.Range("J91").Value = 1627.5
$ws.Range("K91").Value = 2015.5
$ws.Range("L91").Value = 1627.5
$ws.Range("M91").Value = -611.5
$ws.Range("N91").Value = -4435.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1960.0667
$ws.Range("I86").Value = 1761.6154
$ws.Range("J86").Value = 3250
$ws.Range("K86").Value = 1761.6154
$ws.Range("L86").Value = 3250
$ws.Range("M86").Value = -638.6153999999999
$ws.Range("N86").Value = -5496

$ws.Range("H89").Value = 1960.0667
$ws.Range("I89").Value = 1761.6154
$ws.Range("J89").Value = 3250
$ws.Range("K89").Value = 8808.076999999999
$ws.Range("L89").Value = 16250
$ws.Range("M89").Value = -3192.076999999999
$ws.Range("N89").Value = -27482

$ws.Range("H107").Value = 2431
$ws.Range("I107").Value = 2800.3333
$ws.Range("J107").Value = 1600
$ws.Range("K107").Value = 2800.3333
$ws.Range("L107").Value = 1600
$ws.Range("M107").Value = -880.3332999999998
$ws.Range("N107").Value = -5440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2729.6206
$ws.Range("I31").Value = 1956.55
$ws.Range("J31").Value = 4447.5557
$ws.Range("K31").Value = 1956.55
$ws.Range("L31").Value = 4447.5557
$ws.Range("M31").Value = -1661.55
$ws.Range("N31").Value = -5037.5557

$ws.Range("H34").Value = 2729.6206
$ws.Range("I34").Value = 1956.55
$ws.Range("J34").Value = 4447.5557
$ws.Range("K34").Value = 1956.55
$ws.Range("L34").Value = 4447.5557
$ws.Range("M34").Value = -1754.55
$ws.Range("N34").Value = -4851.5557

$ws.Range("H52").Value = 35780
$ws.Range("J52").Value = 35780
$ws.Range("L52").Value = 35780
$ws.Range("N52").Value = -36368

$ws.Range("H119").Value = 51380.5
$ws.Range("J119").Value = 51380.5
$ws.Range("L119").Value = 51380.5
$ws.Range("N119").Value = -61056.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1132.9565
$ws.Range("I17").Value = 313.54544
$ws.Range("J17").Value = 1884.0834
$ws.Range("K17").Value = 940.63632
$ws.Range("L17").Value = 5652.2502
$ws.Range("M17").Value = -771.63632
$ws.Range("N17").Value = -5990.2502

$ws.Range("H34").Value = 398.68182
$ws.Range("I34").Value = 138.06667
$ws.Range("J34").Value = 957.1429000000001
$ws.Range("K34").Value = 414.20001
$ws.Range("L34").Value = 2871.4287
$ws.Range("M34").Value = -330.20001
$ws.Range("N34").Value = -3039.4287

$ws.Range("H39").Value = 1500
$ws.Range("I39").Value = 1500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 4500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -4206
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 2668.1428
$ws.Range("I55").Value = 274
$ws.Range("J55").Value = 3067.1667
$ws.Range("K55").Value = 822
$ws.Range("L55").Value = 9201.500100000001
$ws.Range("M55").Value = -645
$ws.Range("N55").Value = -9555.500100000001

$ws.Range("H110").Value = 3500
$ws.Range("I110").Value = 3500
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 10500
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -6410
$ws.Range("N110").ClearContents()

$ws.Range("H111").Value = 1115.7142
$ws.Range("I111").Value = 1115.7142
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3347.1426
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -280.1425999999997
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 3025
$ws.Range("J112").Value = 3760
$ws.Range("L112").Value = 11280
$ws.Range("N112").Value = -13496

$ws.Range("H120").Value = 10520
$ws.Range("I120").Value = 5153.3335
$ws.Range("J120").Value = 16960
$ws.Range("K120").Value = 15460.0005
$ws.Range("L120").Value = 50880
$ws.Range("M120").Value = -10622.0005
$ws.Range("N120").Value = -60556

$ws.Range("H131").Value = 3787.853
$ws.Range("I131").Value = 655.1429000000001
$ws.Range("J131").Value = 4600.037
$ws.Range("K131").Value = 1965.4287
$ws.Range("L131").Value = 13800.111
$ws.Range("M131").Value = 3074.5713
$ws.Range("N131").Value = -23880.111

$ws.Range("H132").Value = 843226.25
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 1123635
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 10112715
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -10117775

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6391.4165
$ws.Range("I70").Value = 6451.913
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 6451.913
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -6181.913
$ws.Range("N70").Value = -5540

$ws.Range("H73").Value = 6391.4165
$ws.Range("I73").Value = 6451.913
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 6451.913
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -5515.913
$ws.Range("N73").Value = -6872

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4939.2
$ws.Range("I40").Value = 4878.8
$ws.Range("J40").Value = 5060
$ws.Range("K40").Value = 4878.8
$ws.Range("L40").Value = 5060
$ws.Range("M40").Value = -4742.8
$ws.Range("N40").Value = -5332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1887
$ws.Range("N7").ClearContents()

$ws.Range("H15").Value = 10583.333
$ws.Range("J15").Value = 13500
$ws.Range("L15").Value = 13500
$ws.Range("N15").Value = -14076

$ws.Range("H74").Value = 198000
$ws.Range("J74").Value = 198000
$ws.Range("L74").Value = 198000
$ws.Range("N74").Value = -199872

$ws.Range("H77").Value = 198000
$ws.Range("J77").Value = 198000
$ws.Range("L77").Value = 594000
$ws.Range("N77").Value = -603360
